$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from H1 (existing bold/border/center style) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data values for I and J columns, rows 2-57
$iValues = @(7,3,7,8,7,7,6,8,6,8,8,7,7,7,4,7,5,3,9,5,7,7,7,7,6,7,6,7,6,6,7,8,10,7,5,7,7,7,7,8,12,8,8,8,7,9,8,7,6,9,8,9,9,8,8,8)
$jValues = @(7,4,8,8,8,8,6,8,6,8,9,7,7,7,4,7,6,4,9,6,7,7,7,7,6,7,7,8,6,6,7,8,10,7,6,7,7,8,7,8,12,9,8,9,8,9,9,8,7,9,8,9,9,8,8,8)

for ($r = 2; $r -le 57; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
